$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")
$ws.Range("B9").Value = 1
